$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$text = @'
questions = [
    {
        "title": "Due to an oversight, one of your colleagues has referred to a vendor, ANA, as two separate vendors in the QuickBooks Online Vendor Center once spelled correctly (ANA), and once spelled incorrectly (ANNA).How should you clean up the vendor list?",
        "ques_type": 2,
        "options": [
            "Delete the ANNA account.",
            "Merge the ANA and ANNA accounts.",
            "Make ANA the parent vendor.",
            "Make the ANNA account inactive."
        ],
        "score": "Merge the ANA and ANNA accounts."
    },
    {
        "title": "A customer, ABC, has requested a document showing all their unpaid invoices from the last 365 days.How should you generate this in QuickBooks Online?",
        "ques_type": 2,
        "options": [
            "Sales &gt All Sales &gt Filter &gt select Date Last 365 Days &gt Apply.",
            "Sales &gt All Sales &gt New Transaction &gt Time Activity.",
            "Sales &gt Customers &gt select Customer ABC &gt New Transaction &gt select Statement.",
            "Reports &gt Standard &gt Statement of Cash Flows."
        ],
        "score": "Sales &gt Customers &gt select Customer ABC &gt New Transaction &gt select Statement."
    },
    {
        "title": "When processing bills and expenses and making payments in QuickBooks Online, which of the following is correct?",
        "ques_type": 2,
        "options": [
            "The Pay Bills function is used when paying for a purchase immediately.",
            "Creating a bill will record a transaction as an expense and a payment simultaneously.",
            "Bills are for items purchased or services received now but paid for later.",
            "Bills are used when paying for something by credit card."
        ],
        "score": "Bills are for items purchased or services received now but paid for later."
    },
    {
        "title": "When should you process a journal entry in QuickBooks Online?",
        "ques_type": 2,
        "options": [
            "When correcting errors and processing credit memos.",
            "When processing year-end adjustments and recording expenses.",
            "When correcting errors and processing year-end adjustments.",
            "When processing year-end adjustments and creating customer refunds."
        ],
        "score": "When correcting errors and processing year-end adjustments."
    }
]
'@
$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
$ws.Rows(1).AutoFit()
